$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.118.83'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '2.136.48'
$ws.Range('E3').Value = '  -3.64%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.11'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.598'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.88%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '68.55'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -6.35%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.562'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -7.63%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.04'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -10.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0888'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -7.61%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.29'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -7.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0992'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.10%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.55'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -7.27%  '
$ws.Range('D15').Value = '2.458.36'
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.22'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '2.141.35'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.771'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -8.01%  '
$ws.Range('D19').Value = '40.992.51'
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  -8.70%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '68.94'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.30%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.68'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -8.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '223.17'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.41'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -13.63%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.87'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -10.29%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.47'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -11.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.35'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -8.80%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.13'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.49%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '168.56'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.43'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.26%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '30.68'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0745'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -12.61%  '
$ws.Range('E36').Value = '  -5.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0991'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -10.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.04'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.74%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0278'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -8.26%  '
$ws.Range('E40').Value = '  -5.14%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.48'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -17.94%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.21'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -8.10%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '56.61'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -13.31%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.184'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -7.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.13'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.92%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0947'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '96.08'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -8.70%  '
$ws.Range('E48').Value = '  -5.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.09'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -6.60%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.59'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.12'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -12.70%  '
